# Update the "Metadata" sheet (sheet1): bump version/date, replace the
# duplicated "Contact" row with "Jurisdiction", set the Publisher value,
# and delete the now-redundant second "Contact" row.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it's removed entirely and everything below shifts up one row.
$meta.Rows.Item(11).Delete()

# Update the "Elements" sheet (sheet2): the root Extension row's Short /
# Definition columns get the resource-specific text instead of the
# generic Extension placeholders.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "ACA Market Type"
$elements.Range("L2").Value = "Code for Affordable Care Act (ACA) market type of the associated plan"
